$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts existing rows 2-52 down to 3-53)
$ws.Rows("2:2").Insert()

# The inserted row picks up bold-header formatting from the row above;
# clear it and copy the plain date-style formatting used by the other
# data rows (taken from the now-shifted former row 2, i.e. row 3).
$ws.Range("A2:E2").ClearFormats()
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Write the new row-2 values (2007/2008 data point) and the recomputed
# y_0_forecast (C) / y_1_forecast (E) simulated values for every data row.
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 4.930115226412357
$ws.Range("D2").Value = 2008

$ws.Range("A3").Value = 39583
$ws.Range("B3").Value = 2008
$ws.Range("D3").Value = 2009

$ws.Range("A4").Value = 39765
$ws.Range("B4").Value = 2008
$ws.Range("C4").Value = 1.457587285166628
$ws.Range("D4").Value = 2009

$ws.Range("A5").Value = 39948
$ws.Range("B5").Value = 2009
$ws.Range("D5").Value = 2010

$ws.Range("A6").Value = 40130
$ws.Range("B6").Value = 2009
$ws.Range("C6").Value = -0.9140166223623458
$ws.Range("D6").Value = 2010

$ws.Range("A7").Value = 40310
$ws.Range("B7").Value = 2010
$ws.Range("C7").Value = -2.034793027571991
$ws.Range("D7").Value = 2011
$ws.Range("E7").Value = 0.6952816881563351

$ws.Range("A8").Value = 40494
$ws.Range("B8").Value = 2010
$ws.Range("C8").Value = 2.585942866987878
$ws.Range("D8").Value = 2011
$ws.Range("E8").Value = 4.109775046142405

$ws.Range("A9").Value = 40676
$ws.Range("B9").Value = 2011
$ws.Range("C9").Value = 3.880748577052473
$ws.Range("D9").Value = 2012
$ws.Range("E9").Value = 2.936333428994109

$ws.Range("A10").Value = 40862
$ws.Range("B10").Value = 2011
$ws.Range("C10").Value = 4.253963781362402
$ws.Range("D10").Value = 2012
$ws.Range("E10").Value = 2.863367440851095

$ws.Range("A11").Value = 41044
$ws.Range("B11").Value = 2012
$ws.Range("C11").Value = 2.174296999091507
$ws.Range("D11").Value = 2013
$ws.Range("E11").Value = 1.706732094556851

$ws.Range("A12").Value = 41228
$ws.Range("B12").Value = 2012
$ws.Range("C12").Value = 1.752870900283909
$ws.Range("D12").Value = 2013
$ws.Range("E12").Value = 1.520397254708405

$ws.Range("A13").Value = 41409
$ws.Range("B13").Value = 2013
$ws.Range("C13").Value = -1.942476814417471
$ws.Range("D13").Value = 2014
$ws.Range("E13").Value = 1.288975737543607

$ws.Range("A14").Value = 41592
$ws.Range("B14").Value = 2013
$ws.Range("C14").Value = -1.479696720105139
$ws.Range("D14").Value = 2014
$ws.Range("E14").Value = 2.503951807923066

$ws.Range("A15").Value = 41774
$ws.Range("B15").Value = 2014
$ws.Range("C15").Value = 4.926006686093287
$ws.Range("D15").Value = 2015
$ws.Range("E15").Value = 2.643411312704802

$ws.Range("A16").Value = 41957
$ws.Range("B16").Value = 2014
$ws.Range("C16").Value = 3.900127535411246
$ws.Range("D16").Value = 2015
$ws.Range("E16").Value = 1.194160460927884

$ws.Range("A17").Value = 42137
$ws.Range("B17").Value = 2015
$ws.Range("C17").Value = 0.7614971083056643
$ws.Range("D17").Value = 2016
$ws.Range("E17").Value = 2.348955682567344

$ws.Range("A18").Value = 42321
$ws.Range("B18").Value = 2015
$ws.Range("C18").Value = 0.03947433952959933
$ws.Range("D18").Value = 2016
$ws.Range("E18").Value = 1.459149667419779

$ws.Range("A19").Value = 42503
$ws.Range("B19").Value = 2016
$ws.Range("C19").Value = 3.004126378695804
$ws.Range("D19").Value = 2017
$ws.Range("E19").Value = 2.358460376580607

$ws.Range("A20").Value = 42689
$ws.Range("B20").Value = 2016
$ws.Range("C20").Value = 2.192778679161944
$ws.Range("D20").Value = 2017
$ws.Range("E20").Value = 1.586376095629216

$ws.Range("A21").Value = 42867
$ws.Range("B21").Value = 2017
$ws.Range("C21").Value = 2.384525276921168
$ws.Range("D21").Value = 2018
$ws.Range("E21").Value = 2.349880264276161

$ws.Range("A22").Value = 43053
$ws.Range("B22").Value = 2017
$ws.Range("C22").Value = 3.40836448860673
$ws.Range("D22").Value = 2018
$ws.Range("E22").Value = 2.570658574505469

$ws.Range("A23").Value = 43145
$ws.Range("B23").Value = 2018
$ws.Range("C23").Value = 1.512747556564698
$ws.Range("D23").Value = 2019
$ws.Range("E23").Value = 2.136396488383507

$ws.Range("A24").Value = 43235
$ws.Range("B24").Value = 2018
$ws.Range("C24").Value = 2.677874785158663
$ws.Range("D24").Value = 2019
$ws.Range("E24").Value = 2.507403033330702

$ws.Range("A25").Value = 43326
$ws.Range("B25").Value = 2018
$ws.Range("C25").Value = 2.474496385100733
$ws.Range("D25").Value = 2019
$ws.Range("E25").Value = 2.141109033939292

$ws.Range("A26").Value = 43418
$ws.Range("B26").Value = 2018
$ws.Range("C26").Value = 2.799070570134488
$ws.Range("D26").Value = 2019
$ws.Range("E26").Value = 2.479713128614147

$ws.Range("A27").Value = 43510
$ws.Range("B27").Value = 2019
$ws.Range("C27").Value = 2.786825287537487
$ws.Range("D27").Value = 2020
$ws.Range("E27").Value = 1.555206765808892

$ws.Range("A28").Value = 43600
$ws.Range("B28").Value = 2019
$ws.Range("C28").Value = 3.746278935679004
$ws.Range("D28").Value = 2020
$ws.Range("E28").Value = 1.868164064786093

$ws.Range("A29").Value = 43691
$ws.Range("B29").Value = 2019
$ws.Range("C29").Value = 4.066325724929976
$ws.Range("D29").Value = 2020
$ws.Range("E29").Value = 2.042884966056935

$ws.Range("A30").Value = 43783
$ws.Range("B30").Value = 2019
$ws.Range("C30").Value = 4.195393191694419
$ws.Range("D30").Value = 2020
$ws.Range("E30").Value = 2.359935293525561

$ws.Range("A31").Value = 43875
$ws.Range("B31").Value = 2020
$ws.Range("C31").Value = 1.899045195796845
$ws.Range("D31").Value = 2021
$ws.Range("E31").Value = 2.205372086670465

$ws.Range("A32").Value = 43966
$ws.Range("B32").Value = 2020
$ws.Range("C32").Value = 3.278383090085746
$ws.Range("D32").Value = 2021
$ws.Range("E32").Value = 2.401026764575831

$ws.Range("A33").Value = 44068
$ws.Range("B33").Value = 2020
$ws.Range("C33").Value = 2.133862376612439
$ws.Range("D33").Value = 2021
$ws.Range("E33").Value = 1.412539570439053

$ws.Range("A34").Value = 44159
$ws.Range("B34").Value = 2020
$ws.Range("C34").Value = 1.666553973046048
$ws.Range("D34").Value = 2021
$ws.Range("E34").Value = -0.4512719783814068

$ws.Range("A35").Value = 44251
$ws.Range("B35").Value = 2021
$ws.Range("C35").Value = -1.801793552285647
$ws.Range("D35").Value = 2022
$ws.Range("E35").Value = 1.696367270021448

$ws.Range("A36").Value = 44341
$ws.Range("B36").Value = 2021
$ws.Range("C36").Value = 2.877357105987888
$ws.Range("D36").Value = 2022
$ws.Range("E36").Value = 2.770626214993133

$ws.Range("A37").Value = 44432
$ws.Range("B37").Value = 2021
$ws.Range("C37").Value = 1.773820722495745
$ws.Range("D37").Value = 2022
$ws.Range("E37").Value = 2.446369413529137

$ws.Range("A38").Value = 44525
$ws.Range("B38").Value = 2021
$ws.Range("C38").Value = 1.879266440112803
$ws.Range("D38").Value = 2022
$ws.Range("E38").Value = 1.081814991510499

$ws.Range("A39").Value = 44617
$ws.Range("B39").Value = 2022
$ws.Range("C39").Value = -0.1058999733161259
$ws.Range("D39").Value = 2023
$ws.Range("E39").Value = 2.729902493405767

$ws.Range("A40").Value = 44706
$ws.Range("B40").Value = 2022
$ws.Range("C40").Value = -0.4001895765463725
$ws.Range("D40").Value = 2023
$ws.Range("E40").Value = 2.524513388369543

$ws.Range("A41").Value = 44798
$ws.Range("B41").Value = 2022
$ws.Range("C41").Value = -2.404913754290983
$ws.Range("D41").Value = 2023
$ws.Range("E41").Value = 1.366747064788676

$ws.Range("A42").Value = 44890
$ws.Range("B42").Value = 2022
$ws.Range("C42").Value = -2.620683231370946
$ws.Range("D42").Value = 2023
$ws.Range("E42").Value = -1.174318230871441

$ws.Range("A43").Value = 44981
$ws.Range("B43").Value = 2023
$ws.Range("C43").Value = -3.321926039826262
$ws.Range("D43").Value = 2024
$ws.Range("E43").Value = 1.506553781929298

$ws.Range("A44").Value = 45071
$ws.Range("B44").Value = 2023
$ws.Range("C44").Value = -2.321114556364801
$ws.Range("D44").Value = 2024
$ws.Range("E44").Value = 1.637366184014355

$ws.Range("A45").Value = 45163
$ws.Range("B45").Value = 2023
$ws.Range("C45").Value = -2.901570548279864
$ws.Range("D45").Value = 2024
$ws.Range("E45").Value = 1.175840608617551

$ws.Range("A46").Value = 45254
$ws.Range("B46").Value = 2023
$ws.Range("C46").Value = -3.036556262700274
$ws.Range("D46").Value = 2024
$ws.Range("E46").Value = 0.07123445333143685

$ws.Range("A47").Value = 45345
$ws.Range("B47").Value = 2024
$ws.Range("C47").Value = -1.953081240592103
$ws.Range("D47").Value = 2025
$ws.Range("E47").Value = 0.8938306900690307

$ws.Range("A48").Value = 45436
$ws.Range("B48").Value = 2024
$ws.Range("C48").Value = -0.9245282159112467
$ws.Range("D48").Value = 2025
$ws.Range("E48").Value = 0.9277865862836965

$ws.Range("A49").Value = 45534
$ws.Range("B49").Value = 2024
$ws.Range("C49").Value = -2.859191689251428
$ws.Range("D49").Value = 2025
$ws.Range("E49").Value = -0.3916886409131459

$ws.Range("A50").Value = 45618
$ws.Range("B50").Value = 2024
$ws.Range("C50").Value = -2.953443685011514
$ws.Range("D50").Value = 2025
$ws.Range("E50").Value = -1.196842846539037

$ws.Range("A51").Value = 45713
$ws.Range("B51").Value = 2025
$ws.Range("C51").Value = -0.2732731229103447
$ws.Range("D51").Value = 2026
$ws.Range("E51").Value = 0.7812052424394755

$ws.Range("A52").Value = 45800
$ws.Range("B52").Value = 2025
$ws.Range("C52").Value = 0.2384226118222088
$ws.Range("D52").Value = 2026
$ws.Range("E52").Value = 0.8054161303035379

$ws.Range("A53").Value = 45891
$ws.Range("B53").Value = 2025
$ws.Range("C53").Value = -1.131442475565558
$ws.Range("D53").Value = 2026
$ws.Range("E53").Value = -0.4808467302335195
